# Fix Training Data Issue (#48)
# The "Date" column (BF) held the source filename-derived string
# "6-13-2011-12" for every data row. That was one day off because of how
# NBA stats were shown, so replace it with the correct ISO date
# "2012-06-13" for every data row (BF2:BF31). The header in BF1 ("Date")
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2012-06-13"

# Find the used range so this keeps working even if the sheet's extent
# changes; fall back to the known row count (31) otherwise.
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 31) { $lastRow = 31 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("BF$r")
    if ($cell.Value -eq $null -or $cell.Value -eq "") { continue }

    # Writing an ISO "yyyy-mm-dd" looking string directly makes Excel
    # auto-detect it as a date serial and stamp a date NumberFormat on the
    # cell. Temporarily force a Text format so the literal string is
    # stored as-is, then clear the formatting again so the cell ends up
    # with no style applied at all (matching the rest of the column).
    $cell.NumberFormat = "@"
    $cell.Value = $correctDate
    $cell.ClearFormats()
}
